$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data (above the old row 2),
# pushing the existing data rows down from 2-21 to 5-24.
$ws.Rows.Item(2).Resize(3).Insert()
$ws.Rows.Item(2).Resize(3).ClearFormats()

# Fill in the values for the 3 newly inserted rows (new rows 2-4)
$newTopRows = @(
    @(0.0216857157647609, -0.0339030213654041, 0.0201585534960031),
    @(-0.0120645882561802, -0.0719293802976608, 0.0242818929255008),
    @(0.0125227374956011, -0.0106901414692401, -0.0226020142436027)
)

$r = 2
foreach ($row in $newTopRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Append 7 new rows of data at the bottom (new rows 25-31)
$newBottomRows = @(
    @(0.0004581489483825, -0.0694859251379966, -0.005192354787141),
    @(-0.015118914656341, -0.0235183127224445, 0.047036625444889),
    @(0.0102319931611418, 0.0080939643085002, -0.0006108652451075),
    @(0.0326812900602817, 0.0154243474826216, 0.0415388382971286),
    @(0.0032070425804704, 0.0175623763352632, -0.0198531206697225),
    @(0.0097738439217209, -0.0273362193256616, -0.0233655963093042),
    @(0.0073303831741213, 0.0030543261673301, -0.0368046313524246)
)

$r = 25
foreach ($row in $newBottomRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
